$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.629.17"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").Value = "3.375.48"
$ws.Range("E3").Value = "  +4.73%  "
$ws.Range("D5").Value = "'192.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.96%  "
$ws.Range("D6").Value = "'594.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").Value = "'0.423"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "3.962.99"
$ws.Range("E12").Value = "  +4.74%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "'28.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "69.616.31"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").Value = "3.383.86"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "'451.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.86%  "
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "'13.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("E22").Value = "  +3.86%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "3.515.39"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("E27").Value = "  +4.07%  "
$ws.Range("D28").Value = "'9.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D31").Value = "'23.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "'1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("D34").Value = "'7.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("D37").Value = "'164.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.27%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "'27.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "'6.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "2.746.37"
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").Value = "'25.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.45%  "
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'343.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "'40.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'0.0286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("E50").Value = "  +7.59%  "
$ws.Range("E51").Value = "  +4.89%  "

Write-Host "Updated cryptos list"
